$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.41001060599632
$ws.Range("D2").Value = 0.356293015795174
$ws.Range("E2").Value = 0.3926404006151571
$ws.Range("F2").Value = 0.9919449766321051
$ws.Range("G2").Value = 0.4375848817354751
$ws.Range("H2").Value = 0.5941441820509468
$ws.Range("J2").Value = 0.5416419997120556
$ws.Range("L2").Value = 0.5033688759775146
$ws.Range("M2").Value = 0.4028218309640721
$ws.Range("O2").Value = 2.015990485758905
$ws.Range("B3").Value = 1.32008532015783
$ws.Range("D3").Value = 0.3590302003380472
$ws.Range("E3").Value = 0.3883399168390014
$ws.Range("F3").Value = 1.007595196475258
$ws.Range("G3").Value = 0.4408672511645193
$ws.Range("H3").Value = 0.6002901003361245
$ws.Range("J3").Value = 0.526970396728359
$ws.Range("L3").Value = 0.4432435597438769
$ws.Range("M3").Value = 0.3660694387339873
$ws.Range("O3").Value = 2.035515809570953
$ws.Range("B4").Value = 1.264874843027684
$ws.Range("D4").Value = 0.3608602421195002
$ws.Range("E4").Value = 0.3858484676206473
$ws.Range("F4").Value = 1.0179941853293
$ws.Range("G4").Value = 0.4433316687481863
$ws.Range("H4").Value = 0.6044287552210079
$ws.Range("J4").Value = 0.5181354003748453
$ws.Range("L4").Value = 0.406164936221785
$ws.Range("M4").Value = 0.3434325268575762
$ws.Range("O4").Value = 2.049210091926227
$ws.Range("B5").Value = 1.242378602987827
$ws.Range("D5").Value = 0.361643676466926
$ws.Range("E5").Value = 0.3848709440015199
$ws.Range("F5").Value = 1.022430181236185
$ws.Range("G5").Value = 0.4444485617130667
$ws.Range("H5").Value = 0.6062070446800973
$ws.Range("J5").Value = 0.5145792234501414
$ws.Range("L5").Value = 0.3910155843610994
$ws.Range("M5").Value = 0.3341906064422773
$ws.Range("O5").Value = 2.055218825407877
$ws.Range("B6").Value = 1.23864331277349
$ws.Range("D6").Value = 0.3617760440094102
$ws.Range("E6").Value = 0.3847109159477569
$ws.Range("F6").Value = 1.023178745843886
$ws.Range("G6").Value = 0.4446408145625824
$ws.Range("H6").Value = 0.6065078689032646
$ws.Range("J6").Value = 0.5139914060934956
$ws.Range("L6").Value = 0.388497691800751
$ws.Range("M6").Value = 0.3326549721032919
$ws.Range("O6").Value = 2.056242412820794
$ws.Range("B7").Value = 1.264571438877141
$ws.Range("D7").Value = 0.3608706550885898
$ws.Range("E7").Value = 0.3858351311261217
$ws.Range("F7").Value = 1.018053208052066
$ws.Range("G7").Value = 0.4433462759372375
$ws.Range("H7").Value = 0.604452366397048
$ws.Range("J7").Value = 0.5180872610147418
$ws.Range("L7").Value = 0.4059607852964859
$ws.Range("M7").Value = 0.3433079558977639
$ws.Range("O7").Value = 2.049289394965967
$ws.Range("B8").Value = 1.379004381349091
$ws.Range("D8").Value = 0.3572058547485284
$ws.Range("E8").Value = 0.3911268190372823
$ws.Range("F8").Value = 0.9971771231228637
$ws.Range("G8").Value = 0.4386232724819337
$ws.Range("H8").Value = 0.5961875091457642
$ws.Range("J8").Value = 0.5365475557554902
$ws.Range("L8").Value = 0.4826717995838976
$ws.Range("M8").Value = 0.3901646496403472
$ws.Range("O8").Value = 2.022368388762061
$ws.Range("B9").Value = 1.603386346007369
$ws.Range("D9").Value = 0.3511998319389065
$ws.Range("E9").Value = 0.4026753629466668
$ws.Range("F9").Value = 0.962515999530325
$ws.Range("G9").Value = 0.4329377096978533
$ws.Range("H9").Value = 0.5828783051089061
$ws.Range("J9").Value = 0.5741027522332303
$ws.Range("L9").Value = 0.6317810884826827
$ws.Range("M9").Value = 0.4814668553256922
$ws.Range("O9").Value = 1.983142354511159
$ws.Range("B10").Value = 1.768174249557433
$ws.Range("D10").Value = 0.3475005972194225
$ws.Range("E10").Value = 0.4118607411449915
$ws.Range("F10").Value = 0.9408917103092449
$ws.Range("G10").Value = 0.4309597034821024
$ws.Range("H10").Value = 0.5748693695767173
$ws.Range("J10").Value = 0.6024952602652434
$ws.Range("L10").Value = 0.7404842833495877
$ws.Range("M10").Value = 0.5481683406017055
$ws.Range("O10").Value = 1.962639239110217
$ws.Range("B11").Value = 1.843116508522996
$ws.Range("D11").Value = 0.3459712995426187
$ws.Range("E11").Value = 0.4161888328911587
$ws.Range("F11").Value = 0.9318914726628336
$ws.Range("G11").Value = 0.4305414698975767
$ws.Range("H11").Value = 0.5716106484035777
$ws.Range("J11").Value = 0.6155806505448993
$ws.Range("L11").Value = 0.7897441771930289
$ws.Range("M11").Value = 0.5784262536100613
$ws.Range("O11").Value = 1.955127549862368
$ws.Range("B12").Value = 1.871491069821047
$ws.Range("D12").Value = 0.345414161735178
$ws.Range("E12").Value = 0.4178490366848422
$ws.Range("F12").Value = 0.9286039102368733
$ws.Range("G12").Value = 0.4304526619482374
$ws.Range("H12").Value = 0.5704320050735134
$ws.Range("J12").Value = 0.6205596490882499
$ws.Range("L12").Value = 0.808369451671183
$ws.Range("M12").Value = 0.5898714410537167
$ws.Range("O12").Value = 1.952544880838218
$ws.Range("B13").Value = 1.865380321500425
$ws.Range("D13").Value = 0.3455331756016307
$ws.Range("E13").Value = 0.4174905412030938
$ws.Range("F13").Value = 0.9293065756108732
$ws.Range("G13").Value = 0.4304686892267142
$ws.Range("H13").Value = 0.570683383725239
$ws.Range("J13").Value = 0.6194862804598529
$ws.Range("L13").Value = 0.804359444936523
$ws.Range("M13").Value = 0.5874070953557862
$ws.Range("O13").Value = 1.953089446395609
$ws.Range("B14").Value = 1.845450998288072
$ws.Range("D14").Value = 0.3459250236534857
$ws.Range("E14").Value = 0.4163249944927188
$ws.Range("F14").Value = 0.9316185835229902
$ws.Range("G14").Value = 0.4305327675518384
$ws.Range("H14").Value = 0.5715125707904463
$ws.Range("J14").Value = 0.6159898008517928
$ws.Range("L14").Value = 0.7912770646838112
$ws.Range("M14").Value = 0.5793681166409215
$ws.Range("O14").Value = 1.954909819758853
$ws.Range("B15").Value = 1.833243089698783
$ws.Range("D15").Value = 0.3461679006334251
$ws.Range("E15").Value = 0.4156138229687443
$ws.Range("F15").Value = 0.9330504740801402
$ws.Range("G15").Value = 0.4305810865980959
$ws.Range("H15").Value = 0.5720276835911307
$ws.Range("J15").Value = 0.6138511959611606
$ws.Range("L15").Value = 0.783260001160528
$ws.Range("M15").Value = 0.5744423214440104
$ws.Range("O15").Value = 1.956058974814567
$ws.Range("B16").Value = 1.763276054957544
$ws.Range("D16").Value = 0.3476036207585906
$ws.Range("E16").Value = 0.4115808780703532
$ws.Range("F16").Value = 0.9414967610702263
$ws.Range("G16").Value = 0.4309967537792119
$ws.Range("H16").Value = 0.575090078704747
$ws.Range("J16").Value = 0.601643464589074
$ws.Range("L16").Value = 0.7372611215189409
$ws.Range("M16").Value = 0.5461891530718077
$ws.Range("O16").Value = 1.963166743119501
$ws.Range("B17").Value = 1.720347208921226
$ws.Range("D17").Value = 0.3485236331643051
$ws.Range("E17").Value = 0.4091449274090664
$ws.Range("F17").Value = 0.9468928182101166
$ws.Range("G17").Value = 0.4313753217036549
$ws.Range("H17").Value = 0.5770673036331573
$ws.Range("J17").Value = 0.5941974502968463
$ws.Range("L17").Value = 0.7089929187252437
$ws.Range("M17").Value = 0.5288345567608985
$ws.Range("O17").Value = 1.967992643954716
$ws.Range("B18").Value = 1.695653821918256
$ws.Range("D18").Value = 0.3490672542487161
$ws.Range("E18").Value = 0.4077579312897299
$ws.Range("F18").Value = 0.9500752291553241
$ws.Range("G18").Value = 0.4316383728587994
$ws.Range("H18").Value = 0.5782407463520514
$ws.Range("J18").Value = 0.5899306845919909
$ws.Range("L18").Value = 0.6927160017325491
$ws.Range("M18").Value = 0.5188446943731577
$ws.Range("O18").Value = 1.970939189343937
$ws.Range("B19").Value = 1.687292793238328
$ws.Range("D19").Value = 0.3492538005935302
$ws.Range("E19").Value = 0.4072907477819143
$ws.Range("F19").Value = 0.9511662533686334
$ws.Range("G19").Value = 0.4317352095543612
$ws.Range("H19").Value = 0.5786442685944024
$ws.Range("J19").Value = 0.5884887926713702
$ws.Range("L19").Value = 0.68720189474638
$ws.Range("M19").Value = 0.515460950793809
$ws.Range("O19").Value = 1.97196615182682
$ws.Range("B20").Value = 1.72491726018427
$ws.Range("D20").Value = 0.3484242010333745
$ws.Range("E20").Value = 0.4094027813321262
$ws.Range("F20").Value = 0.9463102470212519
$ws.Range("G20").Value = 0.4313303308834833
$ws.Range("H20").Value = 0.5768530781477352
$ws.Range("J20").Value = 0.5949884407176569
$ws.Range("L20").Value = 0.7120039668124036
$ws.Range("M20").Value = 0.5306828112334472
$ws.Range("O20").Value = 1.967461233984096
$ws.Range("B21").Value = 1.851304855992964
$ws.Range("D21").Value = 0.3458093328002363
$ws.Range("E21").Value = 0.4166667692073034
$ws.Range("F21").Value = 0.9309362143949897
$ws.Range("G21").Value = 0.430512055611743
$ws.Range("H21").Value = 0.5712675154130835
$ws.Range("J21").Value = 0.6170161586815368
$ws.Range("L21").Value = 0.7951204556998164
$ws.Range("M21").Value = 0.5817297125467888
$ws.Range("O21").Value = 1.954368019170033
$ws.Range("B22").Value = 1.93387986935835
$ws.Range("D22").Value = 0.344228401332181
$ws.Range("E22").Value = 0.4215379425988601
$ws.Range("F22").Value = 0.921591673175719
$ws.Range("G22").Value = 0.4303829071196077
$ws.Range("H22").Value = 0.5679397568820121
$ws.Range("J22").Value = 0.631551295049178
$ws.Range("L22").Value = 0.8492759802366834
$ws.Range("M22").Value = 0.6150167661819381
$ws.Range("O22").Value = 1.947337448185209
$ws.Range("B23").Value = 1.889810948456386
$ws.Range("D23").Value = 0.3450604914869615
$ws.Range("E23").Value = 0.4189268703279438
$ws.Range("F23").Value = 0.9265145798580647
$ws.Range("G23").Value = 0.4304146152700241
$ws.Range("H23").Value = 0.569686294472973
$ws.Range("J23").Value = 0.623781095901478
$ws.Range("L23").Value = 0.8203876951225197
$ws.Range("M23").Value = 0.5972579067644261
$ws.Range("O23").Value = 1.950949846436856
$ws.Range("B24").Value = 1.722851180450391
$ws.Range("D24").Value = 0.3484691085217619
$ws.Range("E24").Value = 0.4092861636122507
$ws.Range("F24").Value = 0.9465733778285284
$ws.Range("G24").Value = 0.4313505298163989
$ws.Range("H24").Value = 0.576949815144232
$ws.Range("J24").Value = 0.5946307901523937
$ws.Range("L24").Value = 0.7106427502366728
$ws.Range("M24").Value = 0.5298472541954169
$ws.Range("O24").Value = 1.967700948500891
$ws.Range("B25").Value = 1.542692853217659
$ws.Range("D25").Value = 0.3526988862182279
$ws.Range("E25").Value = 0.3994271691524744
$ws.Range("F25").Value = 0.9712192449146855
$ws.Range("G25").Value = 0.4340909700741591
$ws.Range("H25").Value = 0.5861682648570934
$ws.Range("J25").Value = 0.5638006846137955
$ws.Range("L25").Value = 0.591588906741606
$ws.Range("M25").Value = 0.4568320301291706
$ws.Range("O25").Value = 1.992296917565085
